$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values per the published results update
$ws.Range("H2").Value = 47
$ws.Range("I2").Value = 131
$ws.Range("J2").Value = 494
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 136
$ws.Range("M2").Value = 11
$ws.Range("N2").Value = 75
$ws.Range("P2").Value = 0
$ws.Range("R2").Value = 5
$ws.Range("S2").Value = 35
$ws.Range("T2").Value = 90
$ws.Range("U2").Value = 12
$ws.Range("V2").Value = 713
$ws.Range("W2").Value = 2
$ws.Range("X2").Value = 732
$ws.Range("Y2").Value = 1
$ws.Range("Z2").Value = 15
$ws.Range("AA2").Value = 7
